# Insert a new data row at row 18 (pushing all subsequent rows down by one),
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44749
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112022
$ws.Cells.Item(18, 7).Value = "Arveja Verde"
$ws.Cells.Item(18, 8).Value = "Perfection"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 25
$ws.Cells.Item(18, 11).Value = 42000
$ws.Cells.Item(18, 12).Value = 42000
$ws.Cells.Item(18, 13).Value = 42000
$ws.Cells.Item(18, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 1680
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Make sure the new date cell uses the same style/number format as the
# other date cells in column D (style index 2 in the original workbook).
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
